$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.412554860115051
$ws.Range("B1").Value = 1.693804860115051
$ws.Range("C1").Value = 2.270220041275024
$ws.Range("D1").Value = 5.045341491699219
$ws.Range("E1").Value = 2.042982578277588
